$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 465; this shifts the existing rows 465-494 down to 466-495
# and extends the sheet dimension from R494 to R495.
$ws.Rows(465).Insert()

# Populate the newly inserted row with a fresh weekly price record (same
# Mercado/Region/Categoria template as the surrounding rows).
$ws.Cells.Item(465, 1).Value = 10
$ws.Cells.Item(465, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(465, 3).Value = "La Araucanía"
$ws.Cells.Item(465, 4).Value = 44931
$ws.Cells.Item(465, 5).Value = 9
$ws.Cells.Item(465, 6).Value = 100112040
$ws.Cells.Item(465, 7).Value = "Cilantro"
$ws.Cells.Item(465, 8).Value = "Sin especificar"
$ws.Cells.Item(465, 9).Value = "Primera"
$ws.Cells.Item(465, 10).Value = 40
$ws.Cells.Item(465, 11).Value = 10000
$ws.Cells.Item(465, 12).Value = 10000
$ws.Cells.Item(465, 13).Value = 10000
$ws.Cells.Item(465, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(465, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(465, 16).Value = 5000
$ws.Cells.Item(465, 17).Value = 2
$ws.Cells.Item(465, 18).Value = "Hortaliza"
